$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.713081968153176
$ws.Range("B3").Value = 2.66351744507574
$ws.Range("B4").Value = 2.635517553267221
$ws.Range("B5").Value = 2.624716808793551
$ws.Range("B6").Value = 2.622960095037001
$ws.Range("B7").Value = 2.635369426075727
$ws.Range("B8").Value = 2.695485710538321
$ws.Range("B9").Value = 2.83279803290435
$ws.Range("B10").Value = 2.945711141734705
$ws.Range("B11").Value = 2.999733728464207
$ws.Range("B12").Value = 3.020576208424927
$ws.Range("B13").Value = 3.016070220612846
$ws.Range("B14").Value = 3.001440711720647
$ws.Range("B15").Value = 2.992529986297939
$ws.Range("B16").Value = 2.942234384223582
$ws.Range("B17").Value = 2.912062574849131
$ws.Range("B18").Value = 2.894958525305867
$ws.Range("B19").Value = 2.889210228958177
$ws.Range("B20").Value = 2.915248521634339
$ws.Range("B21").Value = 3.005727270943851
$ws.Range("B22").Value = 3.067107696682967
$ws.Range("B23").Value = 3.034141123989343
$ws.Range("B24").Value = 2.913807401158806
$ws.Range("B25").Value = 2.793552801655892
$ws.Range("C2").Value = 0.2595091417840933
$ws.Range("C3").Value = 0.2316076685602297
$ws.Range("C4").Value = 0.2146314260194231
$ws.Range("C5").Value = 0.2077513920593788
$ws.Range("C6").Value = 0.206611224895056
$ws.Range("C7").Value = 0.2145384874322076
$ws.Range("C8").Value = 0.2498557619057351
$ws.Range("C9").Value = 0.3204014771848165
$ws.Range("C10").Value = 0.3731001915127194
$ws.Range("C11").Value = 0.3972814721047371
$ws.Range("C12").Value = 0.4064696986279159
$ws.Range("C13").Value = 0.4044894359780074
$ws.Range("C14").Value = 0.3980367579684412
$ws.Range("C15").Value = 0.3940884227521906
$ws.Range("C16").Value = 0.3715241945693606
$ws.Range("C17").Value = 0.3577360828468272
$ws.Range("C18").Value = 0.3498250569376751
$ws.Range("C19").Value = 0.3471498337301568
$ws.Range("C20").Value = 0.3592018198675646
$ws.Range("C21").Value = 0.3999312056417352
$ws.Range("C22").Value = 0.4267332965055743
$ws.Range("C23").Value = 0.4124113151762003
$ws.Range("C24").Value = 0.3585391106020097
$ws.Range("C25").Value = 0.3011705251186925
$ws.Range("E2").Value = 0.06320847420602149
$ws.Range("E3").Value = 0.05859126288692451
$ws.Range("E4").Value = 0.05574892713754664
$ws.Range("E5").Value = 0.05458861875722931
$ws.Range("E6").Value = 0.05439582222412653
$ws.Range("E7").Value = 0.05573328729077076
$ws.Range("E8").Value = 0.06161787194267987
$ws.Range("E9").Value = 0.07310808462396068
$ws.Range("E10").Value = 0.08153282894105018
$ws.Range("E11").Value = 0.08536462511649034
$ws.Range("E12").Value = 0.08681575278194487
$ws.Range("E13").Value = 0.08650321831080277
$ws.Range("E14").Value = 0.08548400688293611
$ws.Range("E15").Value = 0.0848597308552641
$ws.Range("E16").Value = 0.08128241617390586
$ws.Range("E17").Value = 0.07908782478831
$ws.Range("E18").Value = 0.07782547912387372
$ws.Range("E19").Value = 0.07739805153768486
$ws.Range("E20").Value = 0.07932144904789595
$ws.Range("E21").Value = 0.08578336936479047
$ws.Range("E22").Value = 0.09000734376050445
$ws.Range("E23").Value = 0.08775279341075759
$ws.Range("E24").Value = 0.07921582954735129
$ws.Range("E25").Value = 0.07000363253250796
$ws.Range("F2").Value = 5.664853704764795
$ws.Range("F3").Value = 5.47632377813008
$ws.Range("F4").Value = 5.3615717116881
$ws.Range("F5").Value = 5.315054884482066
$ws.Range("F6").Value = 5.307345417604836
$ws.Range("F7").Value = 5.360943384916226
$ws.Range("F8").Value = 5.599635234695683
$ws.Range("F9").Value = 6.076052143528557
$ws.Range("F10").Value = 6.431697925398481
$ws.Range("F11").Value = 6.594830752483915
$ws.Range("F12").Value = 6.656807885479736
$ws.Range("F13").Value = 6.643450868545301
$ws.Range("F14").Value = 6.599925545692088
$ws.Range("F15").Value = 6.573291624102694
$ws.Range("F16").Value = 6.421064652360428
$ws.Range("F17").Value = 6.328029509899125
$ws.Range("F18").Value = 6.274644501024113
$ws.Range("F19").Value = 6.256590698543988
$ws.Range("F20").Value = 6.337920122632738
$ws.Range("F21").Value = 6.612704430308327
$ws.Range("F22").Value = 6.793474887833042
$ws.Range("F23").Value = 6.696883210859085
$ws.Range("F24").Value = 6.333448259129767
$ws.Range("F25").Value = 5.946220623456753
$ws.Range("G2").Value = 0.00264677777518969
$ws.Range("G3").Value = 0.002653882069416295
$ws.Range("G4").Value = 0.002658467578734359
$ws.Range("G5").Value = 0.002660392618008382
$ws.Range("G6").Value = 0.002660715682523351
$ws.Range("G7").Value = 0.002658493311929803
$ws.Range("G8").Value = 0.002649181092926419
$ws.Range("G9").Value = 0.002632682712789602
$ws.Range("G10").Value = 0.002621621886510845
$ws.Range("G11").Value = 0.002616817277485655
$ws.Range("G12").Value = 0.00261503030697668
$ws.Range("G13").Value = 0.002615413724046826
$ws.Range("G14").Value = 0.00261666961363311
$ws.Range("G15").Value = 0.002617443099401727
$ws.Range("G16").Value = 0.002621940428896714
$ws.Range("G17").Value = 0.002624757386994022
$ws.Range("G18").Value = 0.002626399007088516
$ws.Range("G19").Value = 0.002626958510517332
$ws.Range("G20").Value = 0.002624455305523359
$ws.Range("G21").Value = 0.002616299850106075
$ws.Range("G22").Value = 0.002611158719491279
$ws.Range("G23").Value = 0.00261388542116919
$ws.Range("G24").Value = 0.002624591807762836
$ws.Range("G25").Value = 0.00263695869955427
$ws.Range("J2").Value = 0.4906333479685685
$ws.Range("J3").Value = 0.4682034706495415
$ws.Range("J4").Value = 0.4543669209867431
$ws.Range("J5").Value = 0.4487110548231215
$ws.Range("J6").Value = 0.4477708209486337
$ws.Range("J7").Value = 0.4542907157877778
$ws.Range("J8").Value = 0.4829122246875102
$ws.Range("J9").Value = 0.5385805673337529
$ws.Range("J10").Value = 0.5792780763050018
$ws.Range("J11").Value = 0.5977655727142803
$ws.Range("J12").Value = 0.6047638962533597
$ws.Range("J13").Value = 0.6032567727705498
$ws.Range("J14").Value = 0.5983413717002577
$ws.Range("J15").Value = 0.5953302596576293
$ws.Range("J16").Value = 0.5780694581267767
$ws.Range("J17").Value = 0.5674747394973565
$ws.Range("J18").Value = 0.5613784143289706
$ws.Range("J19").Value = 0.5593138321862057
$ws.Range("J20").Value = 0.5686028175219349
$ws.Range("J21").Value = 0.5997852017000582
$ws.Range("J22").Value = 0.6201507683061607
$ws.Range("J23").Value = 0.6092821396033798
$ws.Range("J24").Value = 0.5680928299883021
$ws.Range("J25").Value = 0.523561918941013
$ws.Range("L2").Value = 0.1592158123296752
$ws.Range("L3").Value = 0.1627636293401729
$ws.Range("L4").Value = 0.1650913798506348
$ws.Range("L5").Value = 0.1660774820554387
$ws.Range("L6").Value = 0.1662434889382425
$ws.Range("L7").Value = 0.1651045268788991
$ws.Range("L8").Value = 0.1604080741018254
$ws.Range("L9").Value = 0.1523852423030156
$ws.Range("L10").Value = 0.147216736590746
$ws.Range("L11").Value = 0.1450236029937209
$ws.Range("L12").Value = 0.1442158999674916
$ws.Range("L13").Value = 0.1443888386957894
$ws.Range("L14").Value = 0.1449566956871742
$ws.Range("L15").Value = 0.1453074941923944
$ws.Range("L16").Value = 0.1473632485789693
$ws.Range("L17").Value = 0.1486649101100035
$ws.Range("L18").Value = 0.1494284635183938
$ws.Range("L19").Value = 0.1496895420243831
$ws.Range("L20").Value = 0.1485248063667726
$ws.Range("L21").Value = 0.1447892833254798
$ws.Range("L22").Value = 0.142480765945912
$ws.Range("L23").Value = 0.1437006869032089
$ws.Range("L24").Value = 0.1485880998951714
$ws.Range("L25").Value = 0.1544283353247309
$ws.Range("M2").Value = 0.4228447748622628
$ws.Range("M3").Value = 0.423761533887145
$ws.Range("M4").Value = 0.4247530577268499
$ws.Range("M5").Value = 0.4252647610724409
$ws.Range("M6").Value = 0.4253562270656133
$ws.Range("M7").Value = 0.424759523034421
$ws.Range("M8").Value = 0.4230718262669377
$ws.Range("M9").Value = 0.4231706375365363
$ws.Range("M10").Value = 0.4253333408225757
$ws.Range("M11").Value = 0.4267740247605545
$ws.Range("M12").Value = 0.4273855069761083
$ws.Range("M13").Value = 0.4272508774285555
$ws.Range("M14").Value = 0.4268230091391771
$ws.Range("M15").Value = 0.4265695200950219
$ws.Range("M16").Value = 0.4252484013404398
$ws.Range("M17").Value = 0.4245551104614051
$ws.Range("M18").Value = 0.4241993350808926
$ws.Range("M19").Value = 0.4240862519164281
$ws.Range("M20").Value = 0.4246244618190858
$ws.Range("M21").Value = 0.4269468932987195
$ws.Range("M22").Value = 0.4288491227746007
$ws.Range("M23").Value = 0.4277986161394267
$ws.Range("M24").Value = 0.4245929747539776
$ws.Range("M25").Value = 0.4227777915320203
